# Mise à jour de l'application
# Adds a new batch of wellness-tracker entries (date 2025-09-30, Excel serial 45930)
# for 12 players at the bottom of the existing table (rows 420-431), extending the
# existing C*D formula down column I, and updates the view selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 420
$lastNewRow  = 431
$templateRow = 419  # last existing data row - used as a style template

# 1) Stamp out 12 new rows by copying the last existing row so that number
#    formats / fonts / alignment (date style, name style, empty-cell style...)
#    are preserved exactly like every other row in the table.
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Range("A${templateRow}:I${templateRow}").Copy($ws.Range("A${r}:I${r}"))
}

# 2) Fill in the actual values for each new row.
#    Columns: A=Date, B=Nom du joueur, C=Volume, D=Intensité, E=Fatigue,
#             F=Douleur, G=Localisation douleur, H=Plaisir, I=Charge (=C*D)

# Row 420 - Kamal Bafounta
$ws.Range("A420").Value = 45930
$ws.Range("B420").Value = "Kamal Bafounta"
$ws.Range("C420").Value = 70
$ws.Range("D420").Value = 6
$ws.Range("E420").Value = 4
$ws.Range("F420").Value = 4
$ws.Range("G411").Copy($ws.Range("G420"))
$ws.Range("G420").Value = "Cheville"
$ws.Range("H420").Value = 6

# Row 421 - Naim Ighbane
$ws.Range("A421").Value = 45930
$ws.Range("B421").Value = "Naim Ighbane"
$ws.Range("C421").Value = 70
$ws.Range("D421").Value = 5
$ws.Range("E421").Value = 5
$ws.Range("F421").Value = 0
$ws.Range("H421").Value = 2

# Row 422 - Omar Benyounes
$ws.Range("A422").Value = 45930
$ws.Range("B422").Value = "Omar Benyounes"
$ws.Range("C422").Value = 70
$ws.Range("D422").Value = 4
$ws.Range("E422").Value = 1
$ws.Range("F422").Value = 0
$ws.Range("H422").Value = 3

# Row 423 - Maé Clavel
$ws.Range("A423").Value = 45930
$ws.Range("B423").Value = "Maé Clavel"
$ws.Range("C423").Value = 70
$ws.Range("D423").Value = 5
$ws.Range("E423").Value = 3
$ws.Range("F423").Value = 0
$ws.Range("H423").Value = 8

# Row 424 - Yoann Martelat
$ws.Range("A424").Value = 45930
$ws.Range("B424").Value = "Yoann Martelat"
$ws.Range("C424").Value = 70
$ws.Range("D424").Value = 4
$ws.Range("E424").Value = 4
$ws.Range("F424").Value = 3
$ws.Range("H424").Value = 7

# Row 425 - Malik Boussaid
$ws.Range("A425").Value = 45930
$ws.Range("B425").Value = "Malik Boussaid"
$ws.Range("C425").Value = 70
$ws.Range("D425").Value = 2
$ws.Range("E425").Value = 0
$ws.Range("F425").Value = 0
$ws.Range("H425").Value = 10

# Row 426 - Hedi Nasri
$ws.Range("A426").Value = 45930
$ws.Range("B426").Value = "Hedi Nasri"
$ws.Range("C426").Value = 70
$ws.Range("D426").Value = 6
$ws.Range("E426").Value = 4
$ws.Range("F426").Value = 2
$ws.Range("G411").Copy($ws.Range("G426"))
$ws.Range("G426").Value = "Adducteur"
$ws.Range("H426").Value = 7

# Row 427 - Ilan Ihaddadene
$ws.Range("A427").Value = 45930
$ws.Range("B427").Value = "Ilan Ihaddadene"
$ws.Range("C427").Value = 70
$ws.Range("D427").Value = 6
$ws.Range("E427").Value = 6
$ws.Range("F427").Value = 0
$ws.Range("H427").Value = 10

# Row 428 - Naim Dhib
$ws.Range("A428").Value = 45930
$ws.Range("B428").Value = "Naim Dhib"
$ws.Range("C428").Value = 70
$ws.Range("D428").Value = 4
$ws.Range("E428").Value = 6
$ws.Range("F428").Value = 0
$ws.Range("H428").Value = 3

# Row 429 - Amir Kherrab
$ws.Range("A429").Value = 45930
$ws.Range("B429").Value = "Amir Kherrab"
$ws.Range("C429").Value = 70
$ws.Range("D429").Value = 4
$ws.Range("E429").Value = 3
$ws.Range("F429").Value = 0
$ws.Range("H429").Value = 10

# Row 430 - Amir Etien
$ws.Range("A430").Value = 45930
$ws.Range("B430").Value = "Amir Etien"
$ws.Range("C430").Value = 70
$ws.Range("D430").Value = 6
$ws.Range("E430").Value = 7
$ws.Range("F430").Value = 5
$ws.Range("G411").Copy($ws.Range("G430"))
$ws.Range("G430").Value = "Ischio adducteur"
$ws.Range("H430").Value = 4

# Row 431 - Ilyes Boughanmi
$ws.Range("A431").Value = 45930
$ws.Range("B431").Value = "Ilyes Boughanmi"
$ws.Range("C431").Value = 70
$ws.Range("D431").Value = 5
$ws.Range("E431").Value = 5
$ws.Range("F431").Value = 0
$ws.Range("H431").Value = 6

# 3) Extend the "Charge" (Volume*Intensité) formula down through the new rows.
$ws.Range("I${firstNewRow}:I${lastNewRow}").Formula = "=C${firstNewRow}*D${firstNewRow}"

# 4) Update the view so the newly added rows are visible/selected, matching
#    the author's last on-screen position.
$ws.Range("K427").Select()
